$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.22"
$ws.Range("E2").Value = "'5.77%"
$ws.Range("D3").Value = "'31.64"
$ws.Range("E3").Value = "'7.59%"
$ws.Range("D4").Value = "'5.208"
$ws.Range("E4").Value = "'2.20%"
$ws.Range("D5").Value = "'0.07277"
$ws.Range("E5").Value = "'7.90%"
$ws.Range("D6").Value = "'7.792"
$ws.Range("E6").Value = "'6.13%"
$ws.Range("D7").Value = "'3.736"
$ws.Range("E7").Value = "'8.46%"
$ws.Range("D8").Value = "'1.469"
$ws.Range("E8").Value = "'6.19%"
$ws.Range("D9").Value = "'0.9056"
$ws.Range("E9").Value = "'-1.32%"
$ws.Range("D10").Value = "'0.01652"
$ws.Range("E10").Value = "'2,459.58%"
$ws.Range("D11").Value = "'0.1668"
$ws.Range("E11").Value = "'4.89%"
$ws.Range("D12").Value = "'0.07480"
$ws.Range("E12").Value = "'9.17%"
$ws.Range("D13").Value = "'0.07908"
$ws.Range("E13").Value = "'3.30%"
$ws.Range("D14").Value = "'0.02976"
$ws.Range("E14").Value = "'1.84%"
$ws.Range("D15").Value = "'0.09909"
$ws.Range("E15").Value = "'10.32%"
$ws.Range("D16").Value = "'0.001516"
$ws.Range("E16").Value = "'-4.43%"
$ws.Range("D17").Value = "'0.04534"
$ws.Range("E17").Value = "'1.37%"
$ws.Range("D18").Value = "'0.006495"
$ws.Range("E18").Value = "'4.21%"
$ws.Range("D19").Value = "'3.466"
$ws.Range("E19").Value = "'0.34%"
$ws.Range("E20").Value = "'-0.13%"
$ws.Range("D21").Value = "'0.3334"
$ws.Range("E21").Value = "'4.24%"
$ws.Range("E22").Value = "'1.90%"
$ws.Range("D23").Value = "'4.283"
$ws.Range("E23").Value = "'4.96%"
$ws.Range("D24").Value = "'0.1628"
$ws.Range("D25").Value = "'0.001221"
$ws.Range("E25").Value = "'2.39%"
$ws.Range("D26").Value = "'0.004418"
$ws.Range("E26").Value = "'6.71%"
$ws.Range("D27").Value = "'0.0001302"
$ws.Range("E27").Value = "'8.63%"
$ws.Range("D28").Value = "'0.0001744"
$ws.Range("E28").Value = "'8.13%"
$ws.Range("D40").Value = "'0.04473"
$ws.Range("E40").Value = "'5.06%"
$ws.Range("D41").Value = "'0.007184"
$ws.Range("E41").Value = "'6.79%"
$ws.Range("D42").Value = "'0.1341"
$ws.Range("E42").Value = "'7.88%"
$ws.Range("D43").Value = "'0.002333"
$ws.Range("E43").Value = "'7.66%"
$ws.Range("D44").Value = "'0.01337"
$ws.Range("E44").Value = "'11.92%"
$ws.Range("D45").Value = "'0.00006078"
$ws.Range("E45").Value = "'6.95%"
$ws.Range("E47").Value = "'7.42%"
